$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 0.5
$ws.Range("C7").Value = 0.5
